# Fill in the December 2023 "Next Statistics" figures (Circulation, ILL Loans,
# ILL Borrows) for every library row, which had previously been left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(71191, 11009, 12393),
    @(38351, 4354, 4376),
    @(121454, 11464, 10633),
    @(2746, 1098, 267),
    @(78663, 13217, 10303),
    @(8235, 2168, 1753),
    @(9415, 1832, 1016),
    @(4517, 684, 445),
    @(454, 417, 32),
    @(3, 0, 0),
    @(1680, 420, 459),
    @(5131, 2056, 1707),
    @(8185, 3064, 1338),
    @(5883, 2229, 980),
    @(3571, 1478, 322),
    @(28387, 4517, 5500),
    @(2425, 1029, 653),
    @(30054, 3820, 5168),
    @(480, 621, 35),
    @(27963, 3908, 4861),
    @(1854, 822, 308),
    @(32674, 4189, 6405),
    @(128621, 12153, 15572),
    @(10045, 3319, 1633),
    @(0, 0, 0),
    @(8722, 1888, 2078),
    @(2322, 745, 487),
    @(23890, 4446, 4294),
    @(732, 269, 363),
    @(4421, 2865, 659),
    @(26321, 5327, 4625),
    @(16733, 4895, 3583),
    @(9245, 1045, 2107),
    @(92867, 9239, 9362),
    @(13798, 4526, 1987),
    @(42175, 3187, 4566),
    @(1564, 1512, 270),
    @(3098, 776, 1105),
    @(4759, 602, 181),
    @(17782, 935, 555),
    @(426, 190, 79),
    @(1441, 152, 134),
    @(1045, 14, 7),
    @(5523, 1527, 739),
    @(20836, 5840, 3647),
    @(50388, 5711, 7191),
    @(23722, 5707, 2134),
    @(18628, 2116, 3051),
    @(52011, 4890, 8046),
    @(7868, 1229, 1930),
    @(21950, 4751, 3971),
    @(3309, 2218, 1238),
    @(3644, 2066, 225),
    @(6914, 1765, 2308),
    @(20879, 7769, 4484),
    @(25453, 1893, 956),
    @(1108116, 170301, 157081)
)

$startRow = 3
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
}

# Move the active selection back to the default (A1) on the frozen pane,
# matching what Excel saves after entering data and pressing Enter/navigating.
$ws.Range("A1").Select()
